# Append a new record (row 28) to the data table on the active sheet.
# Columns: A=Rok (Year), B=Miasto (City), C=Ulica (Street), D=Nazwa obiektu (Object name)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 1992
$ws.Range("C28").Value = "xXX"
$ws.Range("B28").Value = "XXX"
$ws.Range("D28").Value = "XXX"

# Move/update the current selection shown in the saved workbook view.
$ws.Range("E13").Select() | Out-Null
